$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.012.91'
$ws.Range('E2').Value = '  +1.71%  '

# Row 3
$ws.Range('D3').Value = '3.334.71'
$ws.Range('E3').Value = '  +1.64%  '

# Row 4
$ws.Range('E4').Value = '  +0.15%  '

# Row 5
$ws.Range('D5').Value = "'581.19"
$ws.Range('E5').Value = '  +1.97%  '

# Row 6
$ws.Range('D6').Value = "'177.35"
$ws.Range('E6').Value = '  +1.88%  '

# Row 7
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  -0.13%  '

# Row 8
$ws.Range('D8').Value = "'0.590"

# Row 9
$ws.Range('D9').Value = '3.330.20'
$ws.Range('E9').Value = '  +1.62%  '

# Row 10
$ws.Range('E10').Value = '  +6.55%  '

# Row 11
$ws.Range('E11').Value = '  +2.08%  '

# Row 12
$ws.Range('D12').Value = "'47.14"
$ws.Range('E12').Value = '  +4.13%  '

# Row 13
$ws.Range('E13').Value = '  +2.51%  '

# Row 14
$ws.Range('D14').Value = "'687.54"
$ws.Range('E14').Value = '  +0.57%  '

# Row 15
$ws.Range('D15').Value = '3.875.92'
$ws.Range('E15').Value = '  +1.80%  '

# Row 16
$ws.Range('D16').Value = "'8.44"
$ws.Range('E16').Value = '  +2.24%  '

# Row 17
$ws.Range('D17').Value = '67.974.96'
$ws.Range('E17').Value = '  +1.54%  '

# Row 18
$ws.Range('E18').Value = '  -0.26%  '

# Row 19
$ws.Range('D19').Value = '3.322.92'
$ws.Range('E19').Value = '  +1.30%  '

# Row 20
$ws.Range('D20').Value = "'17.46"
$ws.Range('E20').Value = '  +1.26%  '

# Row 21
$ws.Range('D21').Value = "'11.06"
$ws.Range('E21').Value = '  +3.35%  '

# Row 22
$ws.Range('D22').Value = "'0.898"
$ws.Range('E22').Value = '  +1.65%  '

# Row 23
$ws.Range('D23').Value = "'5.39"
$ws.Range('E23').Value = '  +5.04%  '

# Row 24
$ws.Range('D24').Value = "'17.11"
$ws.Range('E24').Value = '  +1.46%  '

# Row 25
$ws.Range('D25').Value = "'99.06"
$ws.Range('E25').Value = '  +0.38%  '

# Row 26
$ws.Range('E26').Value = '  +1.29%  '

# Row 27
$ws.Range('E27').Value = '  +0.84%  '

# Row 28
$ws.Range('D28').Value = "'9.61"
$ws.Range('E28').Value = '  +4.33%  '

# Row 29
$ws.Range('D29').Value = "'33.03"
$ws.Range('E29').Value = '  -0.83%  '

# Row 30
$ws.Range('E30').Value = '  +3.16%  '

# Row 31
$ws.Range('D31').Value = "'7.09"
$ws.Range('E31').Value = '  +6.00%  '

# Row 32
$ws.Range('D32').Value = "'568.98"
$ws.Range('E32').Value = '  -0.84%  '

# Row 33
$ws.Range('D33').Value = "'11.02"
$ws.Range('E33').Value = '  +2.35%  '

# Row 34
$ws.Range('E34').Value = '  +3.24%  '

# Row 35
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = "'57.51"
$ws.Range('E35').Value = '  +3.87%  '

# Row 36
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = '  -0.07%  '

# Row 37
$ws.Range('D37').Value = '3.710.20'
$ws.Range('E37').Value = '  -4.11%  '

# Row 38
$ws.Range('D38').Value = "'3.38"
$ws.Range('E38').Value = '  +2.76%  '

# Row 39
$ws.Range('D39').Value = "'34.57"
$ws.Range('E39').Value = '  +9.44%  '

# Row 40
$ws.Range('E40').Value = '  +3.53%  '

# Row 41
$ws.Range('D41').Value = "'3.20"
$ws.Range('E41').Value = '  +7.11%  '

# Row 42
$ws.Range('D42').Value = "'2.65"
$ws.Range('E42').Value = '  +2.72%  '

# Row 43
$ws.Range('D43').Value = "'3.38"
$ws.Range('E43').Value = '  +0.18%  '

# Row 44
$ws.Range('E44').Value = '  +2.12%  '

# Row 45
$ws.Range('D45').Value = "'0.337"
$ws.Range('E45').Value = '  +3.77%  '

# Row 46
$ws.Range('E46').Value = '  +1.59%  '

# Row 47
$ws.Range('D47').Value = "'2.67"
$ws.Range('E47').Value = '  +6.27%  '

# Row 48
$ws.Range('E48').Value = '  +1.80%  '

# Row 49
$ws.Range('E49').Value = '  -0.31%  '

# Row 50
$ws.Range('D50').Value = "'1.33"
$ws.Range('E50').Value = '  -2.31%  '

# Row 51
$ws.Range('D51').Value = "'130.00"
$ws.Range('E51').Value = '  -0.33%  '
